# Adds trim/adapter configuration columns (F, G, H) to the PAIRED_END sheet:
#   F = "trim" boolean flag (TRUE/FALSE custom format)
#   G = "read_1_adapter" sequence (only when trim = TRUE)
#   H = "read_2_adapter" sequence (only when trim = TRUE)
# Commit message: "functional for ndj analysis with bcftools single calling"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAIRED_END")

$read1Adapter = "AGATCGGAAGAGCGTCGTGTAGGGAAAGAGTGTAGATCTCGGTGGTCGCCGTATCATT"
$read2Adapter = "GATCGGAAGAGCACACGTCTGAACTCCAGTCACGGATGACTATCTCGTATGCCGTCTTCTGCTTG"
$boolFormat = """TRUE"";""TRUE"";""FALSE"""

# --- Header row (row 1): F1 = "trim", G1 = "read_1_adapter", H1 = "read_2_adapter" ---
$ws.Cells.Item(1, 6).Value = "trim"

$ws.Cells.Item(1, 7).Value = "read_1_adapter"
$ws.Cells.Item(1, 7).Font.Bold = $true

$ws.Cells.Item(1, 8).Value = "read_2_adapter"
$ws.Cells.Item(1, 8).Font.Bold = $true

# --- Data rows 2-20 ---
# Rows 2 (w1118), 3 (oregonr) and 20 (hetsub) are not trimmed.
# Rows 4-19 (ndj_01 .. ndj_16) are trimmed, with adapters filled in.
for ($r = 2; $r -le 20; $r++) {
    $trim = ($r -ge 4) -and ($r -le 19)

    $cellF = $ws.Cells.Item($r, 6)
    $cellF.Value = $trim
    $cellF.NumberFormat = $boolFormat

    if ($trim) {
        $ws.Cells.Item($r, 7).Value = $read1Adapter
        $ws.Cells.Item($r, 8).Value = $read2Adapter
    }
}

# --- Column widths for the new / resized columns ---
$ws.Columns.Item(6).ColumnWidth = 7.26
$ws.Columns.Item(7).ColumnWidth = 74.09
$ws.Columns.Item(8).ColumnWidth = 81.88

# --- Selection moves from D25 to B25 ---
$ws.Activate()
$ws.Range("B25").Select()
